$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cumulative-experience values on row 8 (level 1): each value is
# reduced by 200 compared to the previous figures.
$ws.Range("D8").Value = 1034
$ws.Range("E8").Value = 1651
$ws.Range("F8").Value = 2268
$ws.Range("G8").Value = 2885
$ws.Range("H8").Value = 3502

# Move/update the active selection to F24 as left by the editor.
$ws.Range("F24").Select()
